# Rename several table/relation identifiers across the workbook to append
# a "_table" suffix (Golang Gorm/sqlx renderer change).
#
# The same text appears many times throughout the workbook (shared strings
# reused across several sheets/cells), so rather than hard-coding cell
# addresses we scan every used cell on every worksheet and replace exact
# matches of the old identifier with the new one.

$wb = $excel.ActiveWorkbook

$renames = @{
    "modifications"                         = "modifications_table"
    "addresses"                             = "addresses_table"
    "analyticses"                           = "analyticses_table"
    "analyticses_promotion_id"              = "analyticses_table_promotion_id"
    "analytics_promotions"                  = "analytics_promotions_table"
    "analyticses_segment_id"                = "analyticses_table_segment_id"
    "analytics_segments"                    = "analytics_segments_table"
    "analyticses_business_area_id"          = "analyticses_table_business_area_id"
    "analytics_business_areas"              = "analytics_business_areas_table"
    "client_clients"                        = "client_clients_table"
    "client_clients_promotion_id"           = "client_clients_table_promotion_id"
    "client_clients_segment_id"              = "client_clients_table_segment_id"
    "client_clients_business_area_id"       = "client_clients_table_business_area_id"
    "client_client_addresses"               = "client_client_addresses_table"
    "client_client_addresses_client_id_id"  = "client_client_addresses_table_client_id_id"
    "client_client_addresses_client_id"     = "client_client_addresses_table_client_id"
    "order_orders"                          = "order_orders_table"
    "order_orders_client_id"                = "order_orders_table_client_id"
    "order_orders_delivery_address_id"      = "order_orders_table_delivery_address_id"
    "order_orders_promotion_id"             = "order_orders_table_promotion_id"
    "order_orders_segment_id"               = "order_orders_table_segment_id"
    "order_orders_business_area_id"         = "order_orders_table_business_area_id"
    "order_order_lines"                     = "order_order_lines_table"
    "order_order_lines_order_id"            = "order_order_lines_table_order_id"
    "order_order_lines_item_id"             = "order_order_lines_table_item_id"
    "order_items"                           = "order_items_table"
    "order_order_lines_promotion_id"        = "order_order_lines_table_promotion_id"
    "order_order_lines_segment_id"          = "order_order_lines_table_segment_id"
    "order_order_lines_business_area_id"    = "order_order_lines_table_business_area_id"
}

$totalChanges = 0

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $startRow = $used.Row
    $startCol = $used.Column
    $numRows = $used.Rows.Count
    $numCols = $used.Columns.Count
    $endRow = $startRow + $numRows - 1
    $endCol = $startCol + $numCols - 1

    for ($r = $startRow; $r -le $endRow; $r++) {
        for ($c = $startCol; $c -le $endCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($null -ne $val -and $renames.ContainsKey($val)) {
                $cell.Value = $renames[$val]
                $totalChanges++
            }
        }
    }
}

Write-Output "Total cells updated: $totalChanges"
